$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the reference date in A1 (new price list effective date)
$ws.Range("A1").Value = 45436

# Update unit price column (D) values with the new price list figures.
# All values were recalculated with a uniform price-increase factor.
$priceUpdates = @{
    "D25" = 4186.328
    "D26" = 5159.521
    "D27" = 5646.136
    "D28" = 6894.853
    "D33" = 4723.01
    "D34" = 5671.219
    "D35" = 6469.107
    "D36" = 7478.078
    "D37" = 8852.048000000001
    "D38" = 10283.303
    "D39" = 11879.107
    "D40" = 14204.868
    "D45" = 5867.964
    "D46" = 6941.38
    "D47" = 8376.178
    "D48" = 9220.582
    "D49" = 11360.241
    "D50" = 12612.587
    "D51" = 13832.698
    "D52" = 15779.122
    "D53" = 17782.853
    "D58" = 7911.047
    "D59" = 9220.582
    "D60" = 10340.505
    "D61" = 12165.299
    "D62" = 14025.903
    "D63" = 16297.947
    "D64" = 18212.202
    "D65" = 19428.72
    "D66" = 22076.494
    "D67" = 23793.987
    "D68" = 25547.144
    "D69" = 28946.302
    "D74" = 11735.955
    "D75" = 13131.41
    "D76" = 14507.831
    "D77" = 19178.313
    "D78" = 21826.021
    "D79" = 23733.114
    "D80" = 25976.531
    "D81" = 28345.185
    "D82" = 30914.23
    "D83" = 33647.88
    "D84" = 38964.82
    "D85" = 44324.653
    "D90" = 14419.473
    "D91" = 16154.816
    "D92" = 19965.433
    "D93" = 22756.32
    "D94" = 25547.144
    "D95" = 28409.612
    "D96" = 31236.24
    "D97" = 34456.444
    "D98" = 38034.502
    "D99" = 40789.588
    "D100" = 46908.057
    "D101" = 53062.258
    "D102" = 59252.254
    "D103" = 65442.232
    "D108" = 34993.184
    "D109" = 39251.082
    "D110" = 43437.327
    "D111" = 47945.629
    "D112" = 52357.413
    "D113" = 61220.166
    "D114" = 70129.503
    "D115" = 79074.58
    "D116" = 88413.25199999999
    "D117" = 98145.47199999999
    "D118" = 118075.151
    "D119" = 129167.07
    "D120" = 140258.972
    "D121" = 151708.67
    "D122" = 177076.943
}

foreach ($cellRef in $priceUpdates.Keys) {
    $ws.Range($cellRef).Value = $priceUpdates[$cellRef]
}
